$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 12 from row 11 so formatting (style + blank inline-string cells) carries over,
# then overwrite with the new 2021 data.
$ws.Range("A11:U11").Copy($ws.Range("A12"))

$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 7573
$ws.Range("D12").Value = 4875
$ws.Range("E12").Value = 19783
$ws.Range("F12").Value = 11247
$ws.Range("G12").Value = 32095
$ws.Range("H12").Value = 43644
$ws.Range("I12").Value = 3665
$ws.Range("J12").Value = 1967
$ws.Range("K12").Value = 6862
$ws.Range("M12").Value = 85571
$ws.Range("N12").Value = 21710
$ws.Range("O12").Value = 10086
$ws.Range("P12").Value = 2487
$ws.Range("Q12").Value = 1558
$ws.Range("T12").Value = 840
$ws.Range("U12").Value = 176389
